$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Description" column (M) to the building config table.
# Header rows 3 & 4 get the "Description" label (same look as the other
# header cells in column L).
$ws.Range("L3:L4").Copy()
$ws.Range("M3:M4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M3").Value2 = "Description"
$ws.Range("M4").Value2 = "Description"

# Row 5 holds the column "type" markers (int / string / int[]) - the new
# Description column is a string.
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value2 = "string"

# Data rows 6-45 all get the same descriptive text, styled like the
# "Name" column (wrap text, left aligned).
$ws.Range("G6").Copy()
$ws.Range("M6:M45").PasteSpecial(-4122)
$ws.Range("M6:M45").Value2 = "Storage + 100"

$excel.CutCopyMode = 0

# Column M should share the same width/style as columns I:L.
$ws.Range("I1:L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep the sheet's last selection in sync with the edited range.
$ws.Range("N35").Select() | Out-Null
